# Update: Threat Alert Report - 2026-02-04 14:18
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 content updates ---

# Date column A2: 05-FEB-26 -> 28-MAY-26
# (leading apostrophe keeps it literal text instead of being parsed as a date;
#  PasteSpecial(formats) below restores the original cell style afterwards)
$ws.Range("A2").Value = "'28-MAY-26"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Market Threat Airline C2: Nile Air NP-119 -> Nile Air NP-109
$ws.Range("C2").Value = "Nile Air NP-109"

# Fare figures
$ws.Range("D2").Value = 8103
$ws.Range("E2").Value = 12041
$ws.Range("F2").Value = -3938

# IMPACT column J2: LOW THREAT -> MEDIUM THREAT - MONITOR
$ws.Range("J2").Value = "MEDIUM THREAT - MONITOR"

# Re-color the IMPACT cell fill from green to yellow (both fore+back so the
# solid fill reads as a flat FFFFF3CD, matching the "MEDIUM THREAT" styling)
$ws.Range("J2").Interior.Color = 13497343
$ws.Range("J2").Interior.PatternColor = 13497343

# Widen column J (10th column) from 12 to 25 so the longer label fits
$ws.Columns.Item(10).ColumnWidth = 25
